# Tidied up the final model scenarios
#
# On the "spawning-sites" sheet, add a new column L holding the fixed
# particle count per release event (208000 = F*G*H... here simply the
# literal used throughout), for every site row (2-18), a SUM total of
# that column in L20, and a scratch column N with a couple of ad-hoc
# sanity-check figures (N3/N5 formulas, N4 a free-text note).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("spawning-sites")

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 12).Value = 208000
}

$ws.Cells.Item(3, 14).Formula = "=250000/4000"
$ws.Cells.Item(4, 14).Value = "52+8"
$ws.Cells.Item(5, 14).Formula = "=62*4000"

$ws.Cells.Item(20, 12).Formula = "=SUM(L2:L18)"

# Make this the active sheet/cell, as left by the author
$ws.Activate()
$ws.Range("N5").Select() | Out-Null

$wb.Save()
